# "Add enemy dead animation" -> log a new Defect: "Unusual Attack"
# (the enemy damages the player without the player attacking first),
# found in the Forest scene by Fish.

$wb = $excel.ActiveWorkbook
$todo = $wb.Worksheets.Item(1)   # "Todo "
$defect = $wb.Worksheets.Item(2) # "Defect"

# --- Leave a new selection on the "Todo " sheet (no longer the active tab) ---
$todo.Range("B21").Select()

# --- Append the new defect row (row 4) on the "Defect" sheet ---
$defect.Range("A4").Value = "Unusual Attack"
$defect.Range("B4").Value = "Enemy will damage by player without player attack"
$defect.Range("C4").Value = "Forest"
$defect.Range("D4").Value = "Fish"

# Copy the date style from an existing "Create date" cell so the new cell
# reuses the same number format instead of registering a brand new one.
$defect.Range("F3").Copy()
$defect.Range("F4").PasteSpecial(-4122)
$defect.Range("F4").Value = 43104

# --- Make "Defect" the active sheet/tab with its own new selection ---
$defect.Activate()
$defect.Range("G16").Select()
